$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

$ws.Range("D2").Value = "57.923.24"
$ws.Range("E2").Value = "  +2.10%  "
$ws.Range("D3").Value = "3.061.62"
$ws.Range("E3").Value = "  +1.17%  "
$ws.Range("E4").Value = "  -0.03%  "
Set-TextValue $ws.Range("D5") "515.87"
$ws.Range("E5").Value = "  +1.00%  "
Set-TextValue $ws.Range("D6") "141.21"
$ws.Range("E6").Value = "  +0.88%  "
$ws.Range("E7").Value = "  -0.03%  "
Set-TextValue $ws.Range("D8") "0.436"
$ws.Range("E8").Value = "  +0.92%  "
Set-TextValue $ws.Range("D9") "7.30"
$ws.Range("E9").Value = "  +2.25%  "
$ws.Range("E10").Value = "  -0.56%  "
$ws.Range("E11").Value = "  +1.15%  "
$ws.Range("D12").Value = "3.582.33"
$ws.Range("E12").Value = "  +1.04%  "
$ws.Range("E13").Value = "  +2.76%  "
Set-TextValue $ws.Range("D14") "26.29"
$ws.Range("E14").Value = "  +3.48%  "
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("D16").Value = "57.925.22"
$ws.Range("D17").Value = "3.058.69"
$ws.Range("E17").Value = "  +0.98%  "
$ws.Range("E18").Value = "  +2.39%  "
Set-TextValue $ws.Range("D19") "12.80"
$ws.Range("E19").Value = "  -2.56%  "
Set-TextValue $ws.Range("D20") "8.17"
$ws.Range("E20").Value = "  +1.26%  "
Set-TextValue $ws.Range("D21") "330.97"
$ws.Range("E21").Value = "  -0.97%  "
Set-TextValue $ws.Range("D22") "1.00"
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("E23").Value = "  -0.26%  "
Set-TextValue $ws.Range("D24") "65.34"
$ws.Range("E24").Value = "  +0.74%  "
$ws.Range("E25").Value = "  +1.39%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "0.0₃0906"
$ws.Range("E27").Value = "  -2.42%  "
Set-TextValue $ws.Range("D28") "6.47"
$ws.Range("E28").Value = "  +0.94%  "
Set-TextValue $ws.Range("D29") "7.22"
$ws.Range("E29").Value = "  +5.59%  "
$ws.Range("E30").Value = "  +0.56%  "
$ws.Range("E31").Value = "  +2.88%  "
Set-TextValue $ws.Range("D32") "20.59"
$ws.Range("E32").Value = "  +0.69%  "
Set-TextValue $ws.Range("D33") "154.59"
$ws.Range("E33").Value = "  +0.87%  "
Set-TextValue $ws.Range("D34") "4.52"
$ws.Range("E34").Value = "  +0.39%  "
Set-TextValue $ws.Range("D35") "5.98"
$ws.Range("E35").Value = "  +2.74%  "
Set-TextValue $ws.Range("D36") "27.12"
$ws.Range("E36").Value = "  -0.83%  "
Set-TextValue $ws.Range("D37") "1.27"
$ws.Range("E37").Value = "  +3.12%  "
Set-TextValue $ws.Range("D38") "0.0677"
$ws.Range("E38").Value = "  +1.85%  "
$ws.Range("D39").Value = "3.101.82"
$ws.Range("E39").Value = "  +1.23%  "
Set-TextValue $ws.Range("D40") "3.91"
$ws.Range("E40").Value = "  +2.09%  "
$ws.Range("E41").Value = "  +0.40%  "
$ws.Range("E42").Value = "  -0.08%  "
Set-TextValue $ws.Range("D43") "0.657"
$ws.Range("E43").Value = "  -0.27%  "
$ws.Range("D44").Value = "2.301.32"
$ws.Range("E44").Value = "  +4.51%  "
$ws.Range("E45").Value = "  +4.28%  "
Set-TextValue $ws.Range("D46") "1.38"
$ws.Range("E46").Value = "  +1.51%  "
Set-TextValue $ws.Range("D47") "20.77"
$ws.Range("E47").Value = "  +4.61%  "
Set-TextValue $ws.Range("D48") "0.940"
$ws.Range("E48").Value = "  +0.46%  "
$ws.Range("E49").Value = "  +1.28%  "
Set-TextValue $ws.Range("D50") "0.730"
$ws.Range("E50").Value = "  +8.26%  "
$ws.Range("B51").Value = "Bittensor"
$ws.Range("C51").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D51") "253.29"
$ws.Range("E51").Value = "  +9.06%  "
